$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates -------------------------------------------------
# Row 13 ("Grafiek") now gets its definition filled in column B.
$ws.Range("B13").Value = "formules, functies, data, assenstelsel"

# New terms added to column A (rows 14-18), still without a definition
# in column B (matches the source: their B cells stay empty).
$ws.Range("A14").Value = "Tabel"
$ws.Range("A15").Value = "Vloeiendekromme"
$ws.Range("A16").Value = "Scheurlijn"
$ws.Range("A17").Value = "Zaagtand"
$ws.Range("A18").Value = "Kreukellijn"

# --- Formatting updates ------------------------------------------------
# Row 13's term cell now matches the "has a definition" look (same
# yellow fill used by A1:A12) now that B13 is filled in.
$ws.Range("A13").Interior.Color = $ws.Range("A1").Interior.Color

# The whole term column (A1:A50) gets word-wrap turned on.
$ws.Range("A1:A50").WrapText = $true

# --- View/selection state ----------------------------------------------
$ws.Range("A19").Select() | Out-Null
